$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = 3
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5
$ws.Range("BB2").Value = 301
$ws.Range("G5").Value = 1.3
$ws.Range("I5").Value = 13
$ws.Range("N5").Value = 8
$ws.Range("AU5").Value = 13
$ws.Range("AY5").Value = 67
$ws.Range("G8").Value = 1.38
$ws.Range("K8").Value = 2.4
$ws.Range("O8").Value = 1.22
$ws.Range("P8").Value = 4
$ws.Range("Q8").Value = 1.8
$ws.Range("R8").Value = 2
$ws.Range("W8").Value = 7
$ws.Range("Y8").Value = 9
$ws.Range("AK8").Value = 101
$ws.Range("AO8").Value = 6.5
$ws.Range("M11").Value = 1.06
$ws.Range("N11").Value = 10
$ws.Range("X11").Value = 6
$ws.Range("AO11").Value = 7.5
$ws.Range("AQ11").Value = 23
$ws.Range("G14").Value = 6.5
$ws.Range("H14").Value = 4.33
$ws.Range("I14").Value = 1.45
$ws.Range("K14").Value = 2.25
$ws.Range("Q14").Value = 1.98
$ws.Range("R14").Value = 1.88
$ws.Range("AH14").Value = 6
$ws.Range("AR14").Value = 151
$ws.Range("AX14").Value = 7.5

$wb.Save()
